# Add a new time-log entry (row 5) to Sheet1, mirroring the formatting of
# the existing row 4 entry, and update the active selection accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (number formats / styles) of row 4's date & time cells
# down into row 5 so the new row reuses the same style indexes instead of
# creating new ones.
$ws.Range("A4:C4").Copy($ws.Range("A5:C5"))

# Populate the new row's values.
$ws.Range("A5").Value = 45701                      # Date -> 2/13/2025
$ws.Range("B5").Value = 0.0625                      # Start time -> 1:30 AM
$ws.Range("C5").Value = 0.072916666666666671        # End time -> 1:45 AM
$ws.Range("D5").Value = "Worked on Project proposal"

# Match the workbook's saved selection state.
$ws.Range("D5").Select()
